$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 (Fabrica 3): B4, C4, E4 change
$ws.Range("B4").Value = "Lego"
$ws.Range("C4").Value = "EE.UU."
$ws.Range("E4").Value = 10

# Add new row 5 (Fabrica 4)
# Copy formatting of A4 (bold font, border, centered/top alignment) onto A5
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "Fabrica 4:"

$ws.Range("B5").Value = "Bandai Namco"
$ws.Range("C5").Value = "China"
$ws.Range("D5").Value = 9000
$ws.Range("E5").Value = 8
